# Applies the Fenrir_Profits scheduled-runner update: refreshes the
# currentAveragePrice / LevePrice / LeveProfit columns (H-N) for a set of
# leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2000
$ws.Range("I106").Value = 1000
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1000
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -369
$ws.Range("N106").Value = -4262
$ws.Range("H116").Value = 391682.88
$ws.Range("I116").Value = 7789.4443
$ws.Range("J116").Value = 594920.5600000001
$ws.Range("K116").Value = 7789.4443
$ws.Range("L116").Value = 594920.5600000001
$ws.Range("M116").Value = -4347.4443
$ws.Range("N116").Value = -601804.5600000001
$ws.Range("H137").Value = 137881.81
$ws.Range("I137").Value = 189851.22
$ws.Range("J137").Value = 1462.125
$ws.Range("K137").Value = 569553.66
$ws.Range("L137").Value = 4386.375
$ws.Range("M137").Value = -567003.66
$ws.Range("N137").Value = -9486.375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 837.2727
$ws.Range("I74").Value = 624.7646999999999
$ws.Range("J74").Value = 1559.8
$ws.Range("K74").Value = 624.7646999999999
$ws.Range("L74").Value = 1559.8
$ws.Range("M74").Value = 249.2353000000001
$ws.Range("N74").Value = -3307.8
$ws.Range("H77").Value = 837.2727
$ws.Range("I77").Value = 624.7646999999999
$ws.Range("J77").Value = 1559.8
$ws.Range("K77").Value = 3123.8235
$ws.Range("L77").Value = 7799
$ws.Range("M77").Value = 1244.1765
$ws.Range("N77").Value = -16535

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 39295
$ws.Range("J69").Value = 39295
$ws.Range("L69").Value = 39295
$ws.Range("N69").Value = -40917
$ws.Range("H72").Value = 39295
$ws.Range("J72").Value = 39295
$ws.Range("L72").Value = 117885
$ws.Range("N72").Value = -125997

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4937.9727
$ws.Range("I31").Value = 1150.2632
$ws.Range("J31").Value = 9050.343000000001
$ws.Range("K31").Value = 1150.2632
$ws.Range("L31").Value = 9050.343000000001
$ws.Range("M31").Value = -855.2632000000001
$ws.Range("N31").Value = -9640.343000000001
$ws.Range("H34").Value = 4937.9727
$ws.Range("I34").Value = 1150.2632
$ws.Range("J34").Value = 9050.343000000001
$ws.Range("K34").Value = 1150.2632
$ws.Range("L34").Value = 9050.343000000001
$ws.Range("M34").Value = -948.2632000000001
$ws.Range("N34").Value = -9454.343000000001
$ws.Range("H68").Value = 16473.5
$ws.Range("J68").Value = 16473.5
$ws.Range("L68").Value = 16473.5
$ws.Range("N68").Value = -17971.5
$ws.Range("H70").Value = 32000
$ws.Range("J70").Value = 32000
$ws.Range("L70").Value = 32000
$ws.Range("N70").Value = -32630
$ws.Range("H71").Value = 16473.5
$ws.Range("J71").Value = 16473.5
$ws.Range("L71").Value = 49420.5
$ws.Range("N71").Value = -56908.5
$ws.Range("H73").Value = 32000
$ws.Range("J73").Value = 32000
$ws.Range("L73").Value = 32000
$ws.Range("N73").Value = -34184
$ws.Range("H81").Value = 30960
$ws.Range("I81").Value = 27400
$ws.Range("J81").Value = 33333.332
$ws.Range("K81").Value = 27400
$ws.Range("L81").Value = 33333.332
$ws.Range("M81").Value = -26402
$ws.Range("N81").Value = -35329.332
$ws.Range("H84").Value = 30960
$ws.Range("I84").Value = 27400
$ws.Range("J84").Value = 33333.332
$ws.Range("K84").Value = 82200
$ws.Range("L84").Value = 99999.99600000001
$ws.Range("M84").Value = -77208
$ws.Range("N84").Value = -109983.996
$ws.Range("H93").Value = 7799.1665
$ws.Range("I93").Value = 5911.5625
$ws.Range("J93").Value = 22900
$ws.Range("K93").Value = 5911.5625
$ws.Range("L93").Value = 22900
$ws.Range("M93").Value = -4039.5625
$ws.Range("N93").Value = -26644
$ws.Range("H103").Value = 36420
$ws.Range("I103").Value = 4500
$ws.Range("J103").Value = 57700
$ws.Range("K103").Value = 4500
$ws.Range("L103").Value = 57700
$ws.Range("M103").Value = -3328
$ws.Range("N103").Value = -60044

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 29571.5
$ws.Range("J86").Value = 29571.5
$ws.Range("L86").Value = 29571.5
$ws.Range("N86").Value = -31943.5
$ws.Range("H89").Value = 29571.5
$ws.Range("J89").Value = 29571.5
$ws.Range("L89").Value = 88714.5
$ws.Range("N89").Value = -100570.5
$ws.Range("H122").Value = 90915050
$ws.Range("I122").Value = 142859360
$ws.Range("K122").Value = 428578080
$ws.Range("M122").Value = -428575630
$ws.Range("H132").Value = 25670818
$ws.Range("I132").Value = 37076176
$ws.Range("J132").Value = 8761.166999999999
$ws.Range("K132").Value = 111228528
$ws.Range("L132").Value = 26283.501
$ws.Range("M132").Value = -111225998
$ws.Range("N132").Value = -31343.501
$ws.Range("H141").Value = 40665.8
$ws.Range("J141").Value = 40665.8
$ws.Range("L141").Value = 40665.8
$ws.Range("N141").Value = -51025.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 936.25
$ws.Range("I22").Value = 1115.8823
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1115.8823
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -820.8823
$ws.Range("N22").Value = -1090
$ws.Range("H27").Value = 936.25
$ws.Range("I27").Value = 1115.8823
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 1115.8823
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -1008.8823
$ws.Range("N27").Value = -714
$ws.Range("H40").Value = 3464.4
$ws.Range("I40").Value = 3735.739
$ws.Range("J40").Value = 2572.8572
$ws.Range("K40").Value = 3735.739
$ws.Range("L40").Value = 2572.8572
$ws.Range("M40").Value = -3599.739
$ws.Range("N40").Value = -2844.8572
$ws.Range("H132").Value = 4278884
$ws.Range("I132").Value = 6283365
$ws.Range("J132").Value = 2657.6
$ws.Range("K132").Value = 18850095
$ws.Range("L132").Value = 7972.799999999999
$ws.Range("M132").Value = -18847565
$ws.Range("N132").Value = -13032.8
$ws.Range("H136").Value = 5322.7334
$ws.Range("I136").Value = 7124.1
$ws.Range("J136").Value = 1720
$ws.Range("K136").Value = 21372.3
$ws.Range("L136").Value = 5160
$ws.Range("M136").Value = -18822.3
$ws.Range("N136").Value = -10260

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 580283.25
$ws.Range("I136").Value = 743137.9
$ws.Range("K136").Value = 2229413.7
$ws.Range("M136").Value = -2226863.7
$ws.Range("H140").Value = 54686
$ws.Range("J140").Value = 54686
$ws.Range("L140").Value = 54686
$ws.Range("N140").Value = -65046
